# Generate Report for Handoff
# Rotates the "latest handoff" entry: the previously-latest file
# (c7da35dd-...) is pushed down into a newly appended row, and a new
# file (4eedd5b0-...) takes its place as the latest handoff on row 2,
# on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldMd   = "c7da35dd-16c9-43ef-8b41-37ddcd69b272.md"
$newMd   = "4eedd5b0-82e5-4a66-9c8e-9076f352a1a1.md"
$oldHash = "3d619159bdeea1d0e1485b1ebe70d901eb15d6e1"
$newHash = "01a5907ae43871502b07abdeafcbb2940aa548f9"

$oldMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/1ba1b12e138c8659e2771590b82d1757b5312675/e2e/$oldMd"
$newMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/1ba1b12e138c8659e2771590b82d1757b5312675/e2e/$newMd"

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 currently holds the (soon to be previous) latest handoff file;
# push its contents down to a new row 3 before overwriting row 2.
$ws.Range("A3").Value = $oldMd
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-39-13 18:39:45"

$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-40-13 18:40:33"

# Rebuild hyperlinks: row2 -> new file, row3 -> old file.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMd)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldMdUrl, "", "", $oldMd)

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$oldXlf = "$oldMd.$oldHash.zh-cn.xlf"
$newXlf = "$newMd.$newHash.zh-cn.xlf"
$oldXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bed2f8b6cab64774867ecc4c17aed6f8cfd56f3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldXlf"
$newXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bed2f8b6cab64774867ecc4c17aed6f8cfd56f3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlf"

# Push the previous-latest row down to row 3.
$ws.Range("A3").Value = $oldMd
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $oldXlf
$ws.Range("E3").Value = "2016-03-13 18:37:37"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# Overwrite row 2 with the new latest handoff.
$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $newXlf
$ws.Range("E2").Value = "2016-03-13 18:40:30"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMd)
$ws.Hyperlinks.Add($ws.Range("B2"), $newMdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $newXlfUrl, "", "", $newXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldMdUrl, "", "", $oldMd)
$ws.Hyperlinks.Add($ws.Range("B3"), $oldMdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $oldXlfUrl, "", "", $oldXlf)

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$oldXlf = "$oldMd.$oldHash.de-de.xlf"
$newXlf = "$newMd.$newHash.de-de.xlf"
$oldXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6cd059cc6267e9b504603b23051621dab9cfb8c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldXlf"
$newXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6cd059cc6267e9b504603b23051621dab9cfb8c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlf"

# Push the previous-latest row down to row 3.
$ws.Range("A3").Value = $oldMd
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $oldXlf
$ws.Range("E3").Value = "2016-03-13 18:39:45"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# Overwrite row 2 with the new latest handoff.
$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = $newXlf
$ws.Range("E2").Value = "2016-03-13 18:40:33"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMd)
$ws.Hyperlinks.Add($ws.Range("B2"), $newMdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $newXlfUrl, "", "", $newXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldMdUrl, "", "", $oldMd)
$ws.Hyperlinks.Add($ws.Range("B3"), $oldMdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $oldXlfUrl, "", "", $oldXlf)
